$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column C ("Förändrad") bumped by one day for every data row (2-14) ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46060
}

# --- Rows 6-14: Beteckning (A), Datum (B) and Area (G) cycle between rows ---
# New ordering (after the edit) taken row by row.
$rowsData = @(
    @{ Row = 6;  A = "A 2593-2024";   B = 45313.69204861111; G = 2.3 },
    @{ Row = 7;  A = "A 12651-2022";  B = 44641;              G = 3.2 },
    @{ Row = 8;  A = "A 5792-2024";   B = 45335;              G = 5.6 },
    @{ Row = 9;  A = "A 8194-2025";   B = 45708;              G = 1.9 },
    @{ Row = 10; A = "A 13651-2023";  B = 45006;              G = 2.2 },
    @{ Row = 11; A = "A 50997-2025";  B = 45946;              G = 1.5 },
    @{ Row = 12; A = "A 35642-2023";  B = 45147;              G = 1.2 },
    @{ Row = 13; A = "A 28288-2023";  B = 45099.6349537037;   G = 0.5 },
    @{ Row = 14; A = "A 7333-2025";   B = 45703.35899305555;  G = 0.9 }
)

foreach ($entry in $rowsData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 7).Value = $entry.G
}
